# Apply cryptocurrency price/volume updates as described by the commit diff.
# The workbook's worksheet is already open as $excel.ActiveWorkbook.ActiveSheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.482.08'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '2.101.71'
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '''332.50'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = '''1.003'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '''0.5225'
$ws.Range("E7").Value = '  -1.11%  '
$ws.Range("D8").Value = '''0.4480'
$ws.Range("E8").Value = '  +2.18%  '
$ws.Range("D9").Value = '''53.74'
$ws.Range("E9").Value = '  +17.07%  '
$ws.Range("D10").Value = '''0.08919'
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("D12").Value = '''24.35'
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("D13").Value = '2.095.77'
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("D14").Value = '''6.736'
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("D16").Value = '''96.30'
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").Value = '''0.06621'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").Value = '''19.22'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").Value = '''6.285'
$ws.Range("E22").Value = '  -1.18%  '
$ws.Range("D23").Value = '30.532.57'
$ws.Range("D24").Value = '''12.34'
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").Value = '''2.323'
$ws.Range("E25").Value = '  +2.80%  '
$ws.Range("D26").Value = '2.339.02'
$ws.Range("E26").Value = '  -0.79%  '
$ws.Range("D27").Value = '''22.28'
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("D28").Value = '''2.577'
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("D29").Value = '''163.90'
$ws.Range("E29").Value = '  +0.90%  '
$ws.Range("D30").Value = '''132.23'
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").Value = '''1.194'
$ws.Range("E31").Value = '  +1.96%  '
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = '''1.663'
$ws.Range("E33").Value = '  +8.02%  '
$ws.Range("D34").Value = '''6.150'
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").Value = '''3.922'
$ws.Range("E35").Value = '  -3.28%  '
$ws.Range("D36").Value = '''10.42'
$ws.Range("E36").Value = '  +9.52%  '
$ws.Range("D37").Value = '''0.02574'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").Value = '''0.06766'
$ws.Range("E38").Value = '  +0.41%  '
$ws.Range("D39").Value = '''5.482'
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("D40").Value = '''12.74'
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").Value = '''0.6918'
$ws.Range("E42").Value = '  +1.23%  '
$ws.Range("D43").Value = '''1.254'
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '''2.299'
$ws.Range("E45").Value = '  +3.05%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.6361'
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''13.93'
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("D48").Value = '''3.642'
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("D49").Value = '''1.244'
$ws.Range("E49").Value = '  -2.47%  '
$ws.Range("D50").Value = '''1.226'
$ws.Range("E50").Value = '  +5.94%  '
$ws.Range("D51").Value = '''82.33'
$ws.Range("E51").Value = '  -0.11%  '
